$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The "Context" column (previously living on the Terms sheet, column B) is
# being promoted onto the Relations and Rules sheets instead, so that rules
# (and relations) become directly editable per-context in the app.
# ---------------------------------------------------------------------------

# --- Relations sheet: insert a new "Context" column (B) -------------------
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Range("B1").EntireColumn.Insert() | Out-Null
$wsRelations.Range("B1").Value = "relations~"
$wsRelations.Range("B2").Value = "Context"
$wsRelations.Range("B3:B5").Value = "Braga"

# --- Rules sheet: insert a new "Context" column (B) ------------------------
$wsRules = $wb.Worksheets.Item("Rules")
$wsRules.Range("B1").EntireColumn.Insert() | Out-Null
$wsRules.Range("B1").Value = "rules~"
$wsRules.Range("B2").Value = "Context"
$wsRules.Range("B3").Value = "Braga"

# --- Terms sheet: drop the old "Context" column (B) contents --------------
$wsTerms = $wb.Worksheets.Item("Terms")
$wsTerms.Range("B1:B11").ClearContents()

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping to match the edited workbook state.
# Select in the order the edits were made so the last selection wins as the
# workbook's active tab (Relations).
# ---------------------------------------------------------------------------
$wsRules.Range("B3").Select() | Out-Null
$wsTerms.Range("B1:B1048576").Select() | Out-Null
$wsRelations.Range("B3:B5").Select() | Out-Null
